$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the Date column (A2:A10) from text dates stored as shared strings
# into real date values formatted as "d-mmm" (built-in numFmtId 16).
$ws.Range("A2:A11").NumberFormat = "d-mmm"

$ws.Range("A2").Value = 45407
$ws.Range("A3").Value = 45412
$ws.Range("A4").Value = 45414
$ws.Range("A5").Value = 45415
$ws.Range("A6").Value = 45417
$ws.Range("A7").Value = 45418
$ws.Range("A8").Value = 45419
$ws.Range("A9").Value = 45420
$ws.Range("A10").Value = 45421

# New row of data for 10/05/2024
$ws.Range("A11").Value = 45422
$ws.Range("B11").Value = 2.5
$ws.Range("D11").Value = "adding majors and minors to studentgrades, no combined majors yet"

# Update the active selection to reflect where editing left off
$ws.Range("B12").Select()
